# sales_records_template.xlsx
#
# 1. Show user full name rather than login name for all kinds of data:
#    the "manager" and "salesPerson" placeholder columns are switched to
#    the "*FullName" variants.
# 2. (Region/department required-input validation is implemented outside
#    of this spreadsheet template, so there is nothing to change here for
#    that part.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (salesPerson) and column E (manager) hold the JXLS template
# placeholders for row 2. Swap them for the "FullName" variants -
# write F2 before E2 so the two new shared-string entries end up appended
# in salesPersonFullName, managerFullName order.
$ws.Range("F2").Value = '${record.salesPersonFullName}'
$ws.Range("E2").Value = '${record.managerFullName}'

# Update the saved cell selection to match the cell that was last edited.
$ws.Range("E2").Select()
